# 保安员工导入模板 - update sample row data, bold some header runs,
# de-duplicate a style, and tweak the view/column widths.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Replace the sample data row (row 4) with the new employee record.
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "黄春华"
# B4 (性别/gender) stays "男" - unchanged
$ws.Range("C4").Value = "320582197103016717"
$ws.Range("D4").Value = "13915692786"
$ws.Range("E4").Value = "32052016006687"
$ws.Range("F4").Value = "张家港市保安服务有限公司"
$ws.Range("G4").Value = "白鹿小学"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = "汉族"
$ws.Range("J4").Value = ""
$ws.Range("K4").Value = "1971-03-01"
$ws.Range("L4").Value = "江苏省张家港市常阴沙管理区常东社区四组39号"

# ---------------------------------------------------------------------
# 2. Bold the leading run of a few rich-text header / note cells.
# ---------------------------------------------------------------------
$note = $ws.Range("A2").Characters(1, 5)
$note.Font.Bold = $true
$note.Font.Name = "宋体"

$school = $ws.Range("G3").Characters(1, 4)
$school.Font.Bold = $true
$school.Font.Color = 16777215
$school.Font.Name = "宋体-简"

$hometown = $ws.Range("J3").Characters(1, 2)
$hometown.Font.Bold = $true
$hometown.Font.Color = 16777215
$hometown.Font.Name = "宋体"

$birth = $ws.Range("K3").Characters(1, 4)
$birth.Font.Bold = $true
$birth.Font.Color = 16777215
$birth.Font.Name = "宋体"

$addr = $ws.Range("L3").Characters(1, 2)
$addr.Font.Bold = $true
$addr.Font.Color = 16777215
$addr.Font.Name = "宋体-简"

# ---------------------------------------------------------------------
# 3. J2:L2 used a style that duplicated an existing one - repoint them
#    at the equivalent style already used by B2:I2 so the duplicate can
#    drop out of the style table.
# ---------------------------------------------------------------------
$ws.Range("B2").Copy() | Out-Null
$ws.Range("J2:L2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# 4. View tweaks: drop the frozen/scrolled topLeftCell and move the
#    active selection.
# ---------------------------------------------------------------------
$ws.Range("C11").Select() | Out-Null

# ---------------------------------------------------------------------
# 5. Widen columns F, G and L (G is split out of the old F:J block).
# ---------------------------------------------------------------------
$ws.Columns.Item(6).ColumnWidth = 28.93
$ws.Columns.Item(7).ColumnWidth = 31.36
$ws.Columns.Item(12).ColumnWidth = 56.36
